$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: SCD0277 -> SCD0017
$ws.Name = "SCD0017"

# Update TC_ID values in B5:B7 from "DGS-292" to "SCD0017-007"
$ws.Range("B5").Value = "SCD0017-007"
$ws.Range("B6").Value = "SCD0017-007"
$ws.Range("B7").Value = "SCD0017-007"

# Column B widened to fit the new, longer TC_ID text (closest width the
# engine's rounding can reach to the Excel-computed best-fit width)
$ws.Columns("B").ColumnWidth = 11.5

# Move the viewport / selection to the top-left area, selecting B8
$ws.Range("B8").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
